# Update the LOTOMANIA sheet with the latest contest results (concursos 2865-2874)
# These are appended as new rows 359-368 below the existing historical data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner tuple is one lottery draw: Concurso, followed by the 20 drawn numbers (columns A:U)
$newDraws = @(
    ,(2865,2,17,25,34,40,41,46,49,52,55,57,63,71,75,76,80,84,88,97,99)
    ,(2866,7,16,17,18,23,24,38,41,44,45,50,52,56,58,66,75,76,80,90,94)
    ,(2867,4,5,11,22,27,28,34,36,38,44,47,48,51,62,65,66,72,87,89,93)
    ,(2868,3,18,19,24,25,26,28,31,41,48,62,63,68,73,81,85,88,92,94,96)
    ,(2869,1,11,18,22,27,32,36,37,44,45,47,53,66,70,77,84,86,92,93,95)
    ,(2870,4,6,12,13,25,28,29,38,40,43,44,48,53,60,65,67,68,72,84,92)
    ,(2871,1,5,9,10,12,14,20,22,28,30,31,51,55,67,74,80,86,87,93,97)
    ,(2872,1,10,12,14,16,18,21,27,47,53,57,58,61,69,70,76,80,81,89,90)
    ,(2873,2,3,5,10,11,19,23,24,38,39,49,50,57,60,62,74,75,80,86,92)
    ,(2874,1,4,11,15,18,32,36,37,43,44,45,46,48,50,66,76,78,80,81,86)
)

$startRow = 359
$numRows = $newDraws.Count
$numCols = $newDraws[0].Count

$valueArray = New-Object 'object[,]' $numRows,$numCols
for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $valueArray[$r, $c] = $newDraws[$r][$c]
    }
}

$endRow = $startRow + $numRows - 1
$targetRange = $ws.Range("A$startRow`:U$endRow")
$targetRange.Value = $valueArray

# Reflect where the user was working in the sheet after adding the new rows:
# scrolled so column F / row 339 is the top-left visible cell, with B359:U368 selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 339
$ws.Range("B$startRow`:U$endRow").Select() | Out-Null
